$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,4
$data[0,0] = 0.2249999940395355
$data[0,1] = 0.2025251686573029
$data[0,2] = 0.2000000029802322
$data[0,3] = 0.1900897026062012
$data[1,0] = 0.2249999940395355
$data[1,1] = 0.1902805864810944
$data[1,2] = 0.2000000029802322
$data[1,3] = 0.1902458369731903
$data[2,0] = 0.268750011920929
$data[2,1] = 0.1891886591911316
$data[2,2] = 0.2000000029802322
$data[2,3] = 0.1908315420150757
$data[3,0] = 0.2374999970197678
$data[3,1] = 0.1884809583425522
$data[3,2] = 0.2000000029802322
$data[3,3] = 0.1955219358205795
$data[4,0] = 0.2562499940395355
$data[4,1] = 0.1886744201183319
$data[4,2] = 0.25
$data[4,3] = 0.1901074796915054
$data[5,0] = 0.28125
$data[5,1] = 0.1873501390218735
$data[5,2] = 0.2000000029802322
$data[5,3] = 0.1927498877048492
$data[6,0] = 0.2874999940395355
$data[6,1] = 0.1862628012895584
$data[6,2] = 0.2249999940395355
$data[6,3] = 0.1920082867145538
$data[7,0] = 0.331250011920929
$data[7,1] = 0.1839749068021774
$data[7,2] = 0.1500000059604645
$data[7,3] = 0.1943828165531158
$data[8,0] = 0.356249988079071
$data[8,1] = 0.1816280484199524
$data[8,2] = 0.1500000059604645
$data[8,3] = 0.1926939189434052
$data[9,0] = 0.393750011920929
$data[9,1] = 0.1789968907833099
$data[9,2] = 0.1749999970197678
$data[9,3] = 0.1945878565311432
$data[10,0] = 0.418749988079071
$data[10,1] = 0.1739728152751923
$data[10,2] = 0.2249999940395355
$data[10,3] = 0.1964870542287827
$data[11,0] = 0.456250011920929
$data[11,1] = 0.1701681911945343
$data[11,2] = 0.2249999940395355
$data[11,3] = 0.2020240277051926
$data[12,0] = 0.4749999940395355
$data[12,1] = 0.1627120822668076
$data[12,2] = 0.2249999940395355
$data[12,3] = 0.2052769213914871
$data[13,0] = 0.543749988079071
$data[13,1] = 0.1551084220409393
$data[13,2] = 0.2000000029802322
$data[13,3] = 0.2035467624664307
$data[14,0] = 0.6000000238418579
$data[14,1] = 0.1423233896493912
$data[14,2] = 0.1749999970197678
$data[14,3] = 0.2079560458660126
$data[15,0] = 0.6937500238418579
$data[15,1] = 0.1266935467720032
$data[15,2] = 0.125
$data[15,3] = 0.2172607183456421
$data[16,0] = 0.706250011920929
$data[16,1] = 0.1143978387117386
$data[16,2] = 0.2000000029802322
$data[16,3] = 0.2231287062168121
$data[17,0] = 0.8125
$data[17,1] = 0.09102141857147217
$data[17,2] = 0.1749999970197678
$data[17,3] = 0.2483711540699005
$data[18,0] = 0.9125000238418579
$data[18,1] = 0.07408355176448822
$data[18,2] = 0.1500000059604645
$data[18,3] = 0.2414423227310181
$data[19,0] = 0.9375
$data[19,1] = 0.05895372107625008
$data[19,2] = 0.125
$data[19,3] = 0.2618082165718079
$data[20,0] = 0.9937499761581421
$data[20,1] = 0.0436834953725338
$data[20,2] = 0.125
$data[20,3] = 0.2511889934539795
$data[21,0] = 0.9937499761581421
$data[21,1] = 0.03296905383467674
$data[21,2] = 0.125
$data[21,3] = 0.2738993763923645
$data[22,0] = 0.9937499761581421
$data[22,1] = 0.02646046318113804
$data[22,2] = 0.1000000014901161
$data[22,3] = 0.2780892252922058
$data[23,0] = 1
$data[23,1] = 0.01950326189398766
$data[23,2] = 0.1000000014901161
$data[23,3] = 0.2794030010700226
$data[24,0] = 1
$data[24,1] = 0.01473446004092693
$data[24,2] = 0.1000000014901161
$data[24,3] = 0.286578893661499
$data[25,0] = 1
$data[25,1] = 0.01272809971123934
$data[25,2] = 0.1000000014901161
$data[25,3] = 0.2735411524772644
$data[26,0] = 1
$data[26,1] = 0.00965676736086607
$data[26,2] = 0.125
$data[26,3] = 0.2792800068855286
$data[27,0] = 1
$data[27,1] = 0.007384727708995342
$data[27,2] = 0.125
$data[27,3] = 0.28630131483078
$data[28,0] = 1
$data[28,1] = 0.006034146063029766
$data[28,2] = 0.125
$data[28,3] = 0.2861983180046082
$data[29,0] = 1
$data[29,1] = 0.004385470412671566
$data[29,2] = 0.125
$data[29,3] = 0.2859434187412262

$range = $ws.Range("B2:E31")
$range.Value2 = $data
